# Updated cryptos list on Tue Aug 13 19:35:45 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-3 character (U+2083) used in PEPE's tiny price notation (0.0<sub>3</sub>0813).
$sub3 = [char]0x2083
$pepePrice = "0.0{0}0813" -f $sub3

# Helper: write a value into a cell as plain text, even when the text looks
# like a number (e.g. "521.10", "60.384.05"), without permanently changing
# the cell's number format/style. A leading apostrophe forces Excel to treat
# the entry as text; the original style is restored afterwards so the
# generated file does not end up with the cell pointing at a different
# style index than it had originally.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

# Rows where only the Price (D) and Volume(1h) (E) columns changed.
$priceVolumeUpdates = @(
    @(2,  "60.384.05", "  +2.28%  "),
    @(3,  "2.681.19",  "  +1.84%  "),
    @(4,  "0.999",     "  -0.38%  "),
    @(5,  "521.10",    "  +1.30%  "),
    @(6,  "146.21",    "  +1.20%  "),
    @(7,  "0.996",     "  -0.29%  "),
    @(8,  "0.580",     "  +1.67%  "),
    @(9,  "2.689.55",  "  +1.10%  "),
    @(10, "6.47",      "  +2.61%  "),
    @(12, "0.341",     "  +1.01%  "),
    @(14, "3.147.97",  "  +0.78%  "),
    @(15, "60.436.10", "  +2.34%  "),
    @(16, "21.28",     "  +0.80%  "),
    @(17, "2.750.88",  "  +3.39%  "),
    @(18, "0.0000139", "  +1.18%  "),
    @(19, "350.95",    "  +2.01%  "),
    @(20, "4.56",      "  +0.05%  "),
    @(21, "10.56",     "  +1.56%  "),
    @(22, "6.33",      "  +3.68%  "),
    @(23, "1.00",      "  +0.13%  "),
    @(24, "62.91",     "  +2.96%  "),
    @(25, "0.422",     "  +0.24%  "),
    @(26, "0.168",     "  +4.81%  "),
    @(27, "0.990",     "  -0.46%  "),
    @(28, $pepePrice,  "  +0.59%  "),
    @(29, "7.26",      "  +1.90%  "),
    @(30, "6.88",      "  +6.92%  "),
    @(31, "0.998",     "  -0.18%  "),
    @(33, "19.08",     "  +0.85%  "),
    @(34, "148.60",    "  -0.75%  "),
    @(35, "4.32",      "  +6.92%  "),
    @(36, "0.952",     "  -5.72%  "),
    @(37, "1.23",      "  +6.45%  "),
    @(38, "1.57",      "  +11.15%  "),
    @(39, "0.871",     "  +1.85%  "),
    @(40, "36.74",     "  +0.68%  "),
    @(41, "3.70",      "  +0.13%  "),
    @(42, "282.54",    "  +0.64%  "),
    @(43, "0.0990",    "  +0.52%  "),
    @(47, "2.125.39",  "  +7.26%  "),
    @(48, "0.0541",    "  +0.84%  ")
)

foreach ($item in $priceVolumeUpdates) {
    Set-TextValue $item[0] 4 $item[1]
    $ws.Cells.Item($item[0], 5).Value = $item[2]
}

# Rows where only the Volume(1h) (E) column changed.
$volumeOnlyUpdates = @(
    @(11, "  -0.21%  "),
    @(13, "  +1.45%  "),
    @(32, "  +0.93%  "),
    @(44, "  +0.11%  ")
)

foreach ($item in $volumeOnlyUpdates) {
    $ws.Cells.Item($item[0], 5).Value = $item[1]
}

# Rows 45-46 swap: EnergySwap <-> Mantle (now with updated price/volume values).
$ws.Cells.Item(45, 2).Value = "Mantle"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue 45 4 "0.610"
$ws.Cells.Item(45, 5).Value = "  -0.84%  "

$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 46 4 "19.89"
$ws.Cells.Item(46, 5).Value = "  +1.81%  "

# Rows 49-51 rotate: RenderToken -> VeChain -> WhiteBITCoin -> RenderToken.
$ws.Cells.Item(49, 2).Value = "VeChain"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 49 4 "0.0235"
$ws.Cells.Item(49, 5).Value = "  +2.32%  "

$ws.Cells.Item(50, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue 50 4 "10.44"
$ws.Cells.Item(50, 5).Value = "  +1.75%  "

$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 51 4 "4.77"
$ws.Cells.Item(51, 5).Value = "  +3.20%  "
